# Fruta / hortaliza, semanal
# Insert two new weekly data rows into the Guayaba sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at sheet row 5 (pushes former rows 5-18 down to 6-19) ---
$ws.Rows.Item(5).EntireRow.Insert()

$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = "10/4/2021"
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100108
$ws.Cells.Item(5, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(5, 9).Value = 100108001
$ws.Cells.Item(5, 10).Value = "Guayaba"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 160
$ws.Cells.Item(5, 14).Value = 1500
$ws.Cells.Item(5, 15).Value = 1600
$ws.Cells.Item(5, 16).Value = 1550
$ws.Cells.Item(5, 17).Value = "$/kilo (en caja de 10 kilos )"
$ws.Cells.Item(5, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 19).Value = 1550
$ws.Cells.Item(5, 20).Value = 1

# --- Insert new row at sheet row 20 (pushes former rows 19-28, now at 20-29, down to 21-30) ---
$ws.Rows.Item(20).EntireRow.Insert()

$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(20, 4).Value = "3/5/2021"
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100108
$ws.Cells.Item(20, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(20, 9).Value = 100108001
$ws.Cells.Item(20, 10).Value = "Guayaba"
$ws.Cells.Item(20, 11).Value = "Sin especificar"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 100
$ws.Cells.Item(20, 14).Value = 1900
$ws.Cells.Item(20, 15).Value = 2000
$ws.Cells.Item(20, 16).Value = 1950
$ws.Cells.Item(20, 17).Value = "$/kilo (en caja de 10 kilos )"
$ws.Cells.Item(20, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20, 19).Value = 1950
$ws.Cells.Item(20, 20).Value = 1
